$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 250.4
$ws.Range("I55").Value = 158.33333
$ws.Range("J55").Value = 388.5
$ws.Range("K55").Value = 158.33333
$ws.Range("L55").Value = 388.5
$ws.Range("M55").Value = 55.66667000000001
$ws.Range("N55").Value = -816.5
# Row 74
$ws.Range("H74").Value = 19033.334
$ws.Range("I74").Value = 18550
$ws.Range("K74").Value = 18550
$ws.Range("M74").Value = -17614
# Row 76
$ws.Range("H76").Value = 6700.722
$ws.Range("I76").Value = 5566.4
$ws.Range("K76").Value = 5566.4
$ws.Range("M76").Value = -5251.4
# Row 77
$ws.Range("H77").Value = 19033.334
$ws.Range("I77").Value = 18550
$ws.Range("K77").Value = 92750
$ws.Range("M77").Value = -88070
# Row 79
$ws.Range("H79").Value = 6700.722
$ws.Range("I79").Value = 5566.4
$ws.Range("K79").Value = 5566.4
$ws.Range("M79").Value = -4474.4
# Row 100
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2831.3333
$ws.Range("I63").Value = 2831.3333
$ws.Range("K63").Value = 2831.3333
$ws.Range("M63").Value = -2145.3333
# Row 66
$ws.Range("H66").Value = 2831.3333
$ws.Range("I66").Value = 2831.3333
$ws.Range("K66").Value = 14156.6665
$ws.Range("M66").Value = -10724.6665
# Row 88
$ws.Range("H88").Value = 3124
$ws.Range("I88").Value = 1196
$ws.Range("J88").Value = 3766.6667
$ws.Range("K88").Value = 1196
$ws.Range("L88").Value = 3766.6667
$ws.Range("M88").Value = -790
$ws.Range("N88").Value = -4578.6667
# Row 91
$ws.Range("H91").Value = 3124
$ws.Range("I91").Value = 1196
$ws.Range("J91").Value = 3766.6667
$ws.Range("K91").Value = 1196
$ws.Range("L91").Value = 3766.6667
$ws.Range("M91").Value = 208
$ws.Range("N91").Value = -6574.6667

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2095.2058
$ws.Range("I20").Value = 1856.591
$ws.Range("K20").Value = 1856.591
$ws.Range("M20").Value = -1609.591
# Row 86
$ws.Range("H86").Value = 5016.625
$ws.Range("I86").Value = 2991.75
$ws.Range("J86").Value = 7041.5
$ws.Range("K86").Value = 2991.75
$ws.Range("L86").Value = 7041.5
$ws.Range("M86").Value = -1868.75
$ws.Range("N86").Value = -9287.5
# Row 89
$ws.Range("H89").Value = 5016.625
$ws.Range("I89").Value = 2991.75
$ws.Range("J89").Value = 7041.5
$ws.Range("K89").Value = 14958.75
$ws.Range("L89").Value = 35207.5
$ws.Range("M89").Value = -9342.75
$ws.Range("N89").Value = -46439.5
# Row 99
$ws.Range("H99").Value = 1224.9375
$ws.Range("I99").Value = 1299.9286
$ws.Range("J99").Value = 700
$ws.Range("K99").Value = 1299.9286
$ws.Range("L99").Value = 700
$ws.Range("M99").Value = 198.0714
$ws.Range("N99").Value = -3696
# Row 105
$ws.Range("H105").Value = 2357640
$ws.Range("I105").Value = 3080761
$ws.Range("J105").Value = 7497
$ws.Range("K105").Value = 3080761
$ws.Range("L105").Value = 7497
$ws.Range("M105").Value = -3079014
$ws.Range("N105").Value = -10991

$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 109998.6
$ws.Range("J59").Value = 109998.6
$ws.Range("L59").Value = 109998.6
$ws.Range("N59").Value = -112288.6
# Row 62
$ws.Range("H62").Value = 5478.6
$ws.Range("I62").Value = 4939.4
$ws.Range("J62").Value = 6017.8
$ws.Range("K62").Value = 4939.4
$ws.Range("L62").Value = 6017.8
$ws.Range("M62").Value = -4315.4
$ws.Range("N62").Value = -7265.8
# Row 65
$ws.Range("H65").Value = 5478.6
$ws.Range("I65").Value = 4939.4
$ws.Range("J65").Value = 6017.8
$ws.Range("K65").Value = 24697
$ws.Range("L65").Value = 30089
$ws.Range("M65").Value = -21577
$ws.Range("N65").Value = -36329
# Row 68
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101498
# Row 71
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -307488
# Row 99
$ws.Range("H99").Value = 1516.6086
$ws.Range("I99").Value = 1499.2632
$ws.Range("K99").Value = 1499.2632
$ws.Range("M99").Value = -1.263200000000097
# Row 126
$ws.Range("H126").Value = 1516.6086
$ws.Range("I126").Value = 1499.2632
$ws.Range("K126").Value = 4497.7896
$ws.Range("M126").Value = -2027.7896
# Row 134
$ws.Range("H134").Value = 8066213.5
$ws.Range("I134").Value = 8773116
$ws.Range("K134").Value = 26319348
$ws.Range("M134").Value = -26316813

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 4348994
$ws.Range("J113").Value = 7247942
$ws.Range("L113").Value = 21743826
$ws.Range("N113").Value = -21748166

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7029.16
$ws.Range("I70").Value = 6622.222
$ws.Range("K70").Value = 6622.222
$ws.Range("M70").Value = -6352.222
# Row 73
$ws.Range("H73").Value = 7029.16
$ws.Range("I73").Value = 6622.222
$ws.Range("K73").Value = 6622.222
$ws.Range("M73").Value = -5686.222
# Row 80
$ws.Range("H80").Value = 14474
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 14474
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 14474
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -16470
# Row 83
$ws.Range("H83").Value = 14474
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 14474
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 72370
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -82354
# Row 102
$ws.Range("H102").Value = 3208.2693
$ws.Range("I102").Value = 1010.9231
$ws.Range("J102").Value = 5405.615
$ws.Range("K102").Value = 1010.9231
$ws.Range("L102").Value = 5405.615
$ws.Range("M102").Value = 611.0769
$ws.Range("N102").Value = -8649.615
# Row 132
$ws.Range("H132").Value = 8615.799999999999
$ws.Range("I132").Value = 7730.5625
$ws.Range("J132").Value = 10189.556
$ws.Range("K132").Value = 23191.6875
$ws.Range("L132").Value = 30568.668
$ws.Range("M132").Value = -20661.6875
$ws.Range("N132").Value = -35628.66800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 132888.38
$ws.Range("I7").Value = 151457.42
$ws.Range("K7").Value = 151457.42
$ws.Range("M7").Value = -151345.42
# Row 22
$ws.Range("H22").Value = 2387.1765
$ws.Range("I22").Value = 2106.5
$ws.Range("J22").Value = 2788.1428
$ws.Range("K22").Value = 2106.5
$ws.Range("L22").Value = 2788.1428
$ws.Range("M22").Value = -1811.5
$ws.Range("N22").Value = -3378.1428
# Row 27
$ws.Range("H27").Value = 2387.1765
$ws.Range("I27").Value = 2106.5
$ws.Range("J27").Value = 2788.1428
$ws.Range("K27").Value = 2106.5
$ws.Range("L27").Value = 2788.1428
$ws.Range("M27").Value = -1999.5
$ws.Range("N27").Value = -3002.1428
# Row 46
$ws.Range("H46").Value = 3775
$ws.Range("I46").Value = 390.8
$ws.Range("J46").Value = 5655.1113
$ws.Range("K46").Value = 390.8
$ws.Range("L46").Value = 5655.1113
$ws.Range("M46").Value = -202.8
$ws.Range("N46").Value = -6031.1113
# Row 68
$ws.Range("H68").Value = 5684.4287
$ws.Range("I68").Value = 5684.4287
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5684.4287
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4935.4287
$ws.Range("N68").Value = $null
# Row 71
$ws.Range("H71").Value = 5684.4287
$ws.Range("I71").Value = 5684.4287
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 28422.1435
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -24678.1435
$ws.Range("N71").Value = $null
# Row 100
$ws.Range("H100").Value = 2581.1538
$ws.Range("I100").Value = 1961.7778
$ws.Range("K100").Value = 1961.7778
$ws.Range("M100").Value = -1420.7778
# Row 122
$ws.Range("H122").Value = 5338.091
$ws.Range("I122").Value = 4944.5
$ws.Range("K122").Value = 14833.5
$ws.Range("M122").Value = -12383.5
# Row 126
$ws.Range("H126").Value = 132888.38
$ws.Range("I126").Value = 151457.42
$ws.Range("K126").Value = 454372.26
$ws.Range("M126").Value = -451902.26
